$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Created On") to hold the new "Status" field.
$ws.Range("D1").EntireColumn.Insert()

# Header
$ws.Range("D1").Value = "Status"

# Data rows
$ws.Range("D2").Value = "Archived"
$ws.Range("D3").Value = "Blocked"
